# Add new experiment blocks (varying knn # of neighbors) below the existing
# tables on Sheet1: three side-by-side tables at rows 25-29 (7500 / 50 / 1000
# knn) and one more table at rows 31-35 (100 knn).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block 1 (row 25-29): three tables side by side, same layout as the
# existing blocks at rows 1-5 / 7-11 / 13-17 / 19-23.
# ---------------------------------------------------------------------

# Row 25 - section titles
$ws.Range("A25").Value = "Avg MFCC (12 coeff) + Avg Delta + 7500 knn"
$ws.Range("F25").Value = "Avg MFCC (12 coeff) + Avg Delta + 50 knn"
$ws.Range("K25").Value = "Avg MFCC (12 coeff) + Avg Delta + 1000 knn"

# Row 26 - "EER" sub-header
$ws.Range("B26").Value = "EER"
$ws.Range("G26").Value = "EER"
$ws.Range("L26").Value = "EER"

# Row 27 - column headers
$ws.Range("A27").Value = "Train"
$ws.Range("B27").Value = "Test: Read"
$ws.Range("C27").Value = "Test: Phone"
$ws.Range("D27").Value = "Test: Mismatch"
$ws.Range("F27").Value = "Train"
$ws.Range("G27").Value = "Test: Read"
$ws.Range("H27").Value = "Test: Phone"
$ws.Range("I27").Value = "Test: Mismatch"
$ws.Range("K27").Value = "Train"
$ws.Range("L27").Value = "Test: Read"
$ws.Range("M27").Value = "Test: Phone"
$ws.Range("N27").Value = "Test: Mismatch"

# Row 28 - "Read" data row
$ws.Range("A28").Value = "Read"
$ws.Range("B28").Value = 33.577199999999998
$ws.Range("C28").Value = 24.093599999999999
$ws.Range("D28").Value = 48.148099999999999
$ws.Range("F28").Value = "Read"
$ws.Range("G28").Value = 11.4634
$ws.Range("H28").Value = 35
$ws.Range("I28").Value = 35.647300000000001
$ws.Range("K28").Value = "Read"
$ws.Range("L28").Value = 31.4634
$ws.Range("M28").Value = 25
$ws.Range("N28").Value = 41.481499999999997

# Row 29 - "Phone" data row
$ws.Range("A29").Value = "Phone"
$ws.Range("B29").Value = 34.146299999999997
$ws.Range("C29").Value = 24.2105
$ws.Range("D29").Value = 48.148099999999999
$ws.Range("F29").Value = "Phone"
$ws.Range("G29").Value = 33.008099999999999
$ws.Range("H29").Value = 12.1637
$ws.Range("I29").Value = 45.028100000000002
$ws.Range("K29").Value = "Phone"
$ws.Range("L29").Value = 32.195099999999996
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = 45.185200000000002

# Merge the title / sub-header cells, matching the existing blocks' layout
$ws.Range("A25:D25").Merge()
$ws.Range("B26:D26").Merge()
$ws.Range("F25:I25").Merge()
$ws.Range("G26:I26").Merge()
$ws.Range("K25:N25").Merge()
$ws.Range("L26:N26").Merge()

$ws.Range("A25:D25").HorizontalAlignment = -4108
$ws.Range("B26:D26").HorizontalAlignment = -4108
$ws.Range("F25:I25").HorizontalAlignment = -4108
$ws.Range("G26:I26").HorizontalAlignment = -4108
$ws.Range("K25:N25").HorizontalAlignment = -4108
$ws.Range("L26:N26").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Block 2 (row 31-35): single table, "100 knn" variant.
# ---------------------------------------------------------------------

# Row 31 - section title
$ws.Range("A31").Value = "Avg MFCC (12 coeff) + Avg Delta + 100 knn"

# Row 32 - "EER" sub-header
$ws.Range("B32").Value = "EER"

# Row 33 - column headers
$ws.Range("A33").Value = "Train"
$ws.Range("B33").Value = "Test: Read"
$ws.Range("C33").Value = "Test: Phone"
$ws.Range("D33").Value = "Test: Mismatch"

# Row 34 - "Read" data row
$ws.Range("A34").Value = "Read"
$ws.Range("B34").Value = 27.642299999999999
$ws.Range("C34").Value = 28.333300000000001
$ws.Range("D34").Value = 41.481499999999997

# Row 35 - "Phone" data row
$ws.Range("A35").Value = "Phone"
$ws.Range("B35").Value = 31.1111
$ws.Range("C35").Value = 21.403500000000001
$ws.Range("D35").Value = 45.185200000000002

$ws.Range("A31:D31").Merge()
$ws.Range("B32:D32").Merge()

$ws.Range("A31:D31").HorizontalAlignment = -4108
$ws.Range("B32:D32").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Update the view: scrolled down a bit, new selection at H34.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H34").Select()
